# Plan_project.xlsx — re-baseline the "% done" column for the
# "luat xet tuyen => luat" / "api/luathoctap/create" / "Khoa => Mon hoc vuot"
# work items, and scroll the sheet view back up a couple of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scroll the visible pane so row 13 (not row 15) is the top-left row,
# matching topLeftCell A15 -> A13 in the sheetView.
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1

# "% done" (column B) progress updates on the Project1 table rows.
$ws.Range("B9").Value = 1        # Cai dat thuat toan Apriori va Demo mot thuat toan -> done
$ws.Range("B10").Value = 1       # Xay dung trang Admin -> done
$ws.Range("B13").Value = 0.6     # Hoc cong nghe moi -> 60%
$ws.Range("B16").Value = 1       # Them xoa sua lien quan hoc tap -> done
$ws.Range("B17").Value = 1       # Them xoa sua lien quan tuyen sinh -> done
$ws.Range("B19").Value = 0.7     # Luat hoc tap: doi voi mot khoa thi mon nao duoc hoc vuot -> 70%
